# Update "想去人数" (interested-people count) figures in both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet to
# reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览": F3 238 -> 242, F4 861 -> 863 ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 242
$wsExhibit.Range("F4").Value = 863

# --- Sheet "全部类型": F4 238 -> 242, F5 861 -> 863 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 242
$wsAll.Range("F5").Value = 863
